# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Cells in column D that look like plain numbers are written with a leading
# apostrophe so Excel keeps them as text (matching the original inlineStr
# cells) instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.714.74"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "2.078.91"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'233.94"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").Value = "'0.623"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'58.27"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("D10").Value = "'0.0784"
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("D11").Value = "'0.105"
$ws.Range("E11").Value = "  +2.57%  "
$ws.Range("D12").Value = "'14.88"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").Value = "2.385.80"
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").Value = "'20.94"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("D15").Value = "'0.771"
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "2.065.18"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "37.609.51"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").Value = "'6.18"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").Value = "'71.09"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").Value = "0.0₃0832"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").Value = "'227.83"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("D26").Value = "'169.14"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  +2.28%  "
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'19.43"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.40"
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("D32").Value = "'4.66"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "'0.0629"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("D35").Value = "'2.47"
$ws.Range("E35").Value = "  -4.10%  "
$ws.Range("E36").Value = "  +2.84%  "
$ws.Range("E37").Value = "  -4.45%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("D39").Value = "'5.36"
$ws.Range("E39").Value = "  -5.93%  "
$ws.Range("D40").Value = "'0.0978"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("D41").Value = "'97.94"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("D44").Value = "1.450.62"
$ws.Range("E44").Value = "  -1.99%  "
$ws.Range("D45").Value = "'16.49"
$ws.Range("E45").Value = "  +4.77%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'1.16"
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").Value = "'4.26"
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").Value = "'7.37"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("D51").Value = "2.269.79"
$ws.Range("E51").Value = "  -1.97%  "
